# Team Contribution Log -- add Sprint 3 entries for Michael McGregor and
# Ryan Conyac, and move the active selection down to D16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: Michael McGregor's sprint 3 contribution
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = "Michael McGregor"
$ws.Range("D14").Value = "Helped Ryan with the enemy collisions and debugging"

# Row 15: Ryan Conyac's sprint 3 contribution
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = "Ryan Conyac"
$ws.Range("D15").Value = "Edited video for sprint 3"

# Move the selection to where it landed after the edits
$ws.Range("D16").Select()
